$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 634. This shifts the existing rows 634-713
# down to 635-714 (so the old row 713 becomes the new row 714), matching
# the diff's net effect of appending a duplicate of the last row at the
# bottom while a brand new row is inserted at position 634.
$ws.Rows.Item(634).Insert()

# Populate the newly inserted row 634 with its new values.
$ws.Cells.Item(634,1).Value = 5
$ws.Cells.Item(634,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(634,3).Value = "Maule"
$ws.Cells.Item(634,4).Value = 45142
$ws.Cells.Item(634,5).Value = 7
$ws.Cells.Item(634,6).Value = 100112043
$ws.Cells.Item(634,7).Value = "Pepino ensalada"
$ws.Cells.Item(634,8).Value = "Sin especificar"
$ws.Cells.Item(634,9).Value = "Primera"
$ws.Cells.Item(634,10).Value = 300
$ws.Cells.Item(634,11).Value = 10000
$ws.Cells.Item(634,12).Value = 10000
$ws.Cells.Item(634,13).Value = 10000
$ws.Cells.Item(634,14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(634,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(634,16).Value = 167
$ws.Cells.Item(634,17).Value = 60
$ws.Cells.Item(634,18).Value = "Hortaliza"
